$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Grab references to the sheets we need while the old names still apply.
# ---------------------------------------------------------------------------
$wsEurope    = $wb.Worksheets.Item("europe")
$wsEuropeNew = $wb.Worksheets.Item("europe_new")
$wsDocInfo   = $wb.Worksheets.Item("doc_info")
$wsMixMatch  = $wb.Worksheets.Item("mix_match")
$wsCompute   = $wb.Worksheets.Item("compute method")

# ---------------------------------------------------------------------------
# 2. Clone "europe" twice (after doc_info) to build the new "_100" pair.
#    Cloning preserves the row/column/shared-string layout exactly.
# ---------------------------------------------------------------------------
$wsEurope.Copy($null, $wsDocInfo)
$wsLog100 = $wb.Worksheets.Item($wsDocInfo.Index + 1)
$wsLog100.Name = "europe_100_log"

$wsEurope.Copy($null, $wsLog100)
$wsE100 = $wb.Worksheets.Item($wsLog100.Index + 1)
$wsE100.Name = "europe_100"

# ---------------------------------------------------------------------------
# 3. Rename the original pair to their "_500" equivalents.
# ---------------------------------------------------------------------------
$wsEurope.Name    = "europe_500_log"
$wsEuropeNew.Name = "europe_500"

# ---------------------------------------------------------------------------
# 4. Append a brand-new "company" sheet after europe_100 with the lookup
#    column that used to live at the bottom of doc_info.
# ---------------------------------------------------------------------------
$wsCompany = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wsE100)
$wsCompany.Name = "company"
$wsCompany.Range("A1").Value = "company"
$wsCompany.Range("A2").Value = "id"
$wsCompany.Range("A3").Value = "name"
$wsCompany.Range("A4").Value = "url"
$wsCompany.Range("A5").Value = "other_names"
$wsCompany.Range("A6").Value = "other_urls"
$wsCompany.Range("A7").Value = "pay_way"
$wsCompany.Range("A8").Value = "info"

# ---------------------------------------------------------------------------
# 5. Remove the now-duplicated lookup rows (11:18) from doc_info - they live
#    in the "company" sheet now.
# ---------------------------------------------------------------------------
$wsDocInfo.Rows("11:18").Delete()

# ---------------------------------------------------------------------------
# 6. Update view state (selection / scroll / active tab) on every sheet.
# ---------------------------------------------------------------------------

# compute method: just a new active cell.
$wsCompute.Activate()
$wsCompute.Range("B20").Select()

# europe_500_log (was "europe"): drop the scrolled topLeftCell, new cell I32.
$wsEurope.Activate()
$wsEurope.Range("I32").Select()

# europe_500 (was "europe_new"): select range A1:F12.
$wsEuropeNew.Activate()
$wsEuropeNew.Range("A1:F12").Select()

# doc_info: select A11:F18 (anchored at A11).
$wsDocInfo.Activate()
$wsDocInfo.Range("A11:F18").Select()

# europe_100_log: scrolled view with active cell E38.
$wsLog100.Activate()
$wsLog100.Range("E38").Select()

# europe_100: active cell G8.
$wsE100.Activate()
$wsE100.Range("G8").Select()

# company: active cell J17.
$wsCompany.Activate()
$wsCompany.Range("J17").Select()

# mix_match: becomes the active tab, active cell J17.
$wsMixMatch.Activate()
$wsMixMatch.Range("J17").Select()
